$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("area" / " QUALITY" / 1). This shifts row 3
# ("Organisms" / "OR" / 2) up to become the new row 2.
$ws.Rows.Item(2).Delete()
